# Update countries & provincias Spain
# Applies the daily data refresh described by the commit diff:
#  - Update the "Datos actualizados a ..." timestamp banner (A1)
#  - Update case numbers for Estados Unidos (row 4), Alemania (row 7), China (row 9)
#  - Re-order Niger / Kirguistan / Martinica so Martinica now comes first (with fresh
#    stats) followed by Niger and Kirguistan (carrying forward their previous stats)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp banner -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 5 de Abril de 2020 a las 00:22"

# --- Estados Unidos (row 4) -------------------------------------------
$ws.Range("B4").Value = 308124
$ws.Range("C4").Value = 30963
$ws.Range("E4").Value = 285043
$ws.Range("F4").Value = 8113
$ws.Range("G4").Value = 991
$ws.Range("H4").Value = 8395

# --- Alemania (row 7) ---------------------------------------------------
$ws.Range("B7").Value = 96092
$ws.Range("C7").Value = 4933
$ws.Range("E7").Value = 68248
$ws.Range("G7").Value = 169
$ws.Range("H7").Value = 1444

# --- China (row 9) -------------------------------------------------------
$ws.Range("C9").Value = 0
$ws.Range("G9").Value = 0

# --- Niger / Kirguistan / Martinica reorder (rows 113-115) --------------
# New order: Martinica (fresh numbers), Niger (old row-113 numbers),
# Kirguistan (old row-114 numbers).
$ws.Range("A113").Value = "Martinica"
$ws.Range("B113").Value = 145
$ws.Range("C113").Value = 2
$ws.Range("D113").Value = 27
$ws.Range("E113").Value = 115
$ws.Range("F113").Value = 22
$ws.Range("G113").Value = 0
$ws.Range("H113").Value = 3

$ws.Range("A114").Value = "Niger"
$ws.Range("B114").Value = 144
$ws.Range("C114").Value = 24
$ws.Range("D114").Value = 0
$ws.Range("E114").Value = 136
$ws.Range("F114").Value = 0
$ws.Range("G114").Value = 3
$ws.Range("H114").Value = 8

$ws.Range("A115").Value = "Kirguistan"
$ws.Range("B115").Value = 144
$ws.Range("C115").Value = 14
$ws.Range("D115").Value = 9
$ws.Range("E115").Value = 134
$ws.Range("F115").Value = 5
$ws.Range("G115").Value = 0
$ws.Range("H115").Value = 1
